$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.31316573748933
$ws.Range("C2").Value = 11.84010561515647
$ws.Range("D2").Value = 14.15956757559038
$ws.Range("E2").Value = 14.84811228905033
$ws.Range("G2").Value = 54.50571722401614
$ws.Range("H2").Value = 20.32345201087863
$ws.Range("J2").Value = 8.804551075236112
$ws.Range("M2").Value = 21.10196533253012
$ws.Range("N2").Value = 20.09585692038839
$ws.Range("B3").Value = 18.90903804794742
$ws.Range("C3").Value = 11.47483974471798
$ws.Range("D3").Value = 14.15133265898079
$ws.Range("E3").Value = 14.86717736098309
$ws.Range("G3").Value = 54.16403796566748
$ws.Range("H3").Value = 20.32839887147685
$ws.Range("J3").Value = 8.822916133225236
$ws.Range("M3").Value = 20.99173617891381
$ws.Range("N3").Value = 20.16664817796953
$ws.Range("B4").Value = 18.66227009808558
$ws.Range("C4").Value = 11.24815731486175
$ws.Range("D4").Value = 14.14918911161825
$ws.Range("E4").Value = 14.8813916364922
$ws.Range("G4").Value = 53.97074229940368
$ws.Range("H4").Value = 20.33605387164364
$ws.Range("J4").Value = 8.834873604584125
$ws.Range("M4").Value = 20.92882995850882
$ws.Range("N4").Value = 20.21211858886305
$ws.Range("B5").Value = 18.56220773770978
$ws.Range("C5").Value = 11.15533633640147
$ws.Range("D5").Value = 14.14904841333506
$ws.Range("E5").Value = 14.88781401456687
$ws.Range("G5").Value = 53.89617619942163
$ws.Range("H5").Value = 20.34033170834338
$ws.Range("J5").Value = 8.839918138362195
$ws.Range("M5").Value = 20.90441495526168
$ws.Range("N5").Value = 20.23115327545823
$ws.Range("B6").Value = 18.54562703936355
$ws.Range("C6").Value = 11.13990139318545
$ws.Range("D6").Value = 14.14906931373177
$ws.Range("E6").Value = 14.88891846784225
$ws.Range("G6").Value = 53.88404988847753
$ws.Range("H6").Value = 20.34111190990903
$ws.Range("J6").Value = 8.840766167240719
$ws.Range("M6").Value = 20.90043505734542
$ws.Range("N6").Value = 20.23434451233698
$ws.Range("B7").Value = 18.66091840304684
$ws.Range("C7").Value = 11.24690708684073
$ws.Range("D7").Value = 14.1491842467796
$ws.Range("E7").Value = 14.88147570142258
$ws.Range("G7").Value = 53.96971958860559
$ws.Range("H7").Value = 20.33610687809975
$ws.Range("J7").Value = 8.834940940847352
$ws.Range("M7").Value = 20.92849572652784
$ws.Range("N7").Value = 20.21237325064232
$ws.Range("B8").Value = 19.1736299415763
$ws.Range("C8").Value = 11.71475417408169
$ws.Range("D8").Value = 14.15612399183474
$ws.Range("E8").Value = 14.85416490597774
$ws.Range("G8").Value = 54.38451785003357
$ws.Range("H8").Value = 20.32419779804294
$ws.Range("J8").Value = 8.810742230712004
$ws.Range("M8").Value = 21.06297928459298
$ws.Range("N8").Value = 20.11985043755843
$ws.Range("B9").Value = 20.18331408745182
$ws.Range("C9").Value = 12.60662158141576
$ws.Range("D9").Value = 14.19281117883172
$ws.Range("E9").Value = 14.82054179504506
$ws.Range("G9").Value = 55.32603149663916
$ws.Range("H9").Value = 20.3375846967809
$ws.Range("J9").Value = 8.768673251251952
$ws.Range("M9").Value = 21.36366725189469
$ws.Range("N9").Value = 19.95426377959314
$ws.Range("B10").Value = 20.91920561054616
$ws.Range("C10").Value = 13.23839670334807
$ws.Range("D10").Value = 14.23376725079765
$ws.Range("E10").Value = 14.80802921406976
$ws.Range("G10").Value = 56.09160239746194
$ws.Range("H10").Value = 20.36992500405676
$ws.Range("J10").Value = 8.741018974867218
$ws.Range("M10").Value = 21.6057748460663
$ws.Range("N10").Value = 19.84219420587848
$ws.Range("B11").Value = 21.25100957635809
$ws.Range("C11").Value = 13.51927105140103
$ws.Range("D11").Value = 14.25541637038682
$ws.Range("E11").Value = 14.8049896943582
$ws.Range("G11").Value = 56.45487179497734
$ws.Range("H11").Value = 20.38953455579231
$ws.Range("J11").Value = 8.72913877794209
$ws.Range("M11").Value = 21.72020040904216
$ws.Range("N11").Value = 19.7932766000578
$ws.Range("B12").Value = 21.3760991504071
$ws.Range("C12").Value = 13.62458833265829
$ws.Range("D12").Value = 14.26404554361297
$ws.Range("E12").Value = 14.80422033228817
$ws.Range("G12").Value = 56.59449430633205
$ws.Range("H12").Value = 20.39766415239693
$ws.Range("J12").Value = 8.724740235084665
$ws.Range("M12").Value = 21.76411853529686
$ws.Range("N12").Value = 19.77504837707132
$ws.Range("B13").Value = 21.34918568250884
$ws.Range("C13").Value = 13.60195435938506
$ws.Range("D13").Value = 14.26216797282305
$ws.Range("E13").Value = 14.80436905290547
$ws.Range("G13").Value = 56.56433407518547
$ws.Range("H13").Value = 20.39588200628419
$ws.Range("J13").Value = 8.725683088729912
$ws.Range("M13").Value = 21.75463432133018
$ws.Range("N13").Value = 19.77896100876686
$ws.Range("B14").Value = 21.26131267067123
$ws.Range("C14").Value = 13.52795705218955
$ws.Range("D14").Value = 14.25611767531911
$ws.Range("E14").Value = 14.80491874919074
$ws.Range("G14").Value = 56.46631790867666
$ws.Range("H14").Value = 20.39018928037355
$ws.Range("J14").Value = 8.728774900755553
$ws.Range("M14").Value = 21.72380197338807
$ws.Range("N14").Value = 19.79177102973953
$ws.Range("B15").Value = 21.20741153760321
$ws.Range("C15").Value = 13.48249268914167
$ws.Range("D15").Value = 14.25246774560252
$ws.Range("E15").Value = 14.80530515835909
$ws.Range("G15").Value = 56.40654550309267
$ws.Range("H15").Value = 20.38679396463734
$ws.Range("J15").Value = 8.730681765667333
$ws.Range("M15").Value = 21.70499191441043
$ws.Range("N15").Value = 19.79965603464628
$ws.Range("B16").Value = 20.89745046135583
$ws.Range("C16").Value = 13.21990035233387
$ws.Range("D16").Value = 14.23241292517808
$ws.Range("E16").Value = 14.80828125040092
$ws.Range("G16").Value = 56.06815562720565
$ws.Range("H16").Value = 20.36874200947646
$ws.Range("J16").Value = 8.741809421664719
$ws.Range("M16").Value = 21.59838071417136
$ws.Range("N16").Value = 19.8454325297246
$ws.Range("B17").Value = 20.70644552114938
$ws.Range("C17").Value = 13.05705897400666
$ws.Range("D17").Value = 14.2208810285101
$ws.Range("E17").Value = 14.81078655368603
$ws.Range("G17").Value = 55.86434014097112
$ws.Range("H17").Value = 20.35892207457302
$ws.Range("J17").Value = 8.748814838420358
$ws.Range("M17").Value = 21.53405645462157
$ws.Range("N17").Value = 19.87404276415118
$ws.Range("B18").Value = 20.59631364640904
$ws.Range("C18").Value = 12.96278851851377
$ws.Range("D18").Value = 14.21453245458524
$ws.Range("E18").Value = 14.81247720177139
$ws.Range("G18").Value = 55.74852961142332
$ws.Range("H18").Value = 20.35373497646311
$ws.Range("J18").Value = 8.752910068423217
$ws.Range("M18").Value = 21.49746442259966
$ws.Range("N18").Value = 19.89069288422364
$ws.Range("B19").Value = 20.55898251853962
$ws.Range("C19").Value = 12.93076906950054
$ws.Range("D19").Value = 14.2124318397546
$ws.Range("E19").Value = 14.81309249743855
$ws.Range("G19").Value = 55.70956472615236
$ws.Range("H19").Value = 20.35205790489879
$ws.Range("J19").Value = 8.754307973835489
$ws.Range("M19").Value = 21.48514551926068
$ws.Range("N19").Value = 19.89636372270421
$ws.Range("B20").Value = 20.72680734152866
$ws.Range("C20").Value = 13.07445748444475
$ws.Range("D20").Value = 14.22207921726601
$ws.Range("E20").Value = 14.81049401864508
$ws.Range("G20").Value = 55.88589050546314
$ws.Range("H20").Value = 20.35991970546106
$ws.Range("J20").Value = 8.748062282485288
$ws.Range("M20").Value = 21.5408621137611
$ws.Range("N20").Value = 19.87097705628012
$ws.Range("B21").Value = 21.28713924570543
$ws.Range("C21").Value = 13.54972097001401
$ws.Range("D21").Value = 14.25788311885918
$ws.Range("E21").Value = 14.80474693155815
$ws.Range("G21").Value = 56.49505251979332
$ws.Range("H21").Value = 20.39184227523148
$ws.Range("J21").Value = 8.727864043657977
$ws.Range("M21").Value = 21.73284247070866
$ws.Range("N21").Value = 19.78800039485108
$ws.Range("B22").Value = 21.6500440332453
$ws.Range("C22").Value = 13.85420799056181
$ws.Range("D22").Value = 14.28379430476428
$ws.Range("E22").Value = 14.80321535252554
$ws.Range("G22").Value = 56.90513547812174
$ws.Range("H22").Value = 20.41680764876748
$ws.Range("J22").Value = 8.715247355990476
$ws.Range("M22").Value = 21.86172373550119
$ws.Range("N22").Value = 19.73549417188006
$ws.Range("B23").Value = 21.45669888128376
$ws.Range("C23").Value = 13.6922901096693
$ws.Range("D23").Value = 14.26973629614114
$ws.Range("E23").Value = 14.80382921329195
$ws.Range("G23").Value = 56.6852055816717
$ws.Range("H23").Value = 20.40310811813855
$ws.Range("J23").Value = 8.721927815507195
$ws.Range("M23").Value = 21.79263518540316
$ws.Range("N23").Value = 19.7633602929294
$ws.Range("B24").Value = 20.71760275098489
$ws.Range("C24").Value = 13.06659363676676
$ws.Range("D24").Value = 14.22153663983335
$ws.Range("E24").Value = 14.81062549407636
$ws.Range("G24").Value = 55.87614332444241
$ws.Range("H24").Value = 20.35946724834688
$ws.Range("J24").Value = 8.74840230205157
$ws.Range("M24").Value = 21.53778406190896
$ws.Range("N24").Value = 19.87236243444254
$ws.Range("B25").Value = 19.91061297480015
$ws.Range("C25").Value = 12.36894485635601
$ws.Range("D25").Value = 14.18042080282629
$ws.Range("E25").Value = 14.82749925703606
$ws.Range("G25").Value = 55.05802511047554
$ws.Range("H25").Value = 20.33001637607367
$ws.Range("J25").Value = 8.779480558742049
$ws.Range("M25").Value = 21.27849821838105
$ws.Range("N25").Value = 19.99737039747051
